$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "309.02"
Set-TextCell $ws.Range("E2") "-2.75%"
Set-TextCell $ws.Range("D3") "37.67"
Set-TextCell $ws.Range("E3") "-5.29%"
Set-TextCell $ws.Range("D4") "5.106"
Set-TextCell $ws.Range("E4") "-0.76%"
Set-TextCell $ws.Range("D5") "0.07849"
Set-TextCell $ws.Range("E5") "-4.42%"
Set-TextCell $ws.Range("D6") "1.972"
Set-TextCell $ws.Range("E6") "-1.93%"
Set-TextCell $ws.Range("D7") "4.363"
Set-TextCell $ws.Range("E7") "1.89%"
Set-TextCell $ws.Range("D8") "8.295"
Set-TextCell $ws.Range("E8") "-0.03%"
Set-TextCell $ws.Range("D9") "3.154"
Set-TextCell $ws.Range("E9") "-1.61%"
Set-TextCell $ws.Range("D10") "0.9272"
Set-TextCell $ws.Range("E10") "-0.61%"
Set-TextCell $ws.Range("D11") "0.1351"
Set-TextCell $ws.Range("E11") "-5.30%"
Set-TextCell $ws.Range("D12") "0.1986"
Set-TextCell $ws.Range("E12") "0.25%"
Set-TextCell $ws.Range("D13") "0.08957"
Set-TextCell $ws.Range("E13") "-1.30%"
Set-TextCell $ws.Range("D14") "0.03436"
Set-TextCell $ws.Range("E14") "-2.87%"
Set-TextCell $ws.Range("D15") "0.09728"
Set-TextCell $ws.Range("E15") "-0.76%"
Set-TextCell $ws.Range("D16") "0.001395"
Set-TextCell $ws.Range("E16") "-0.49%"
Set-TextCell $ws.Range("D17") "0.005900"
Set-TextCell $ws.Range("E17") "-7.87%"
Set-TextCell $ws.Range("E18") "1,777.94%"
Set-TextCell $ws.Range("D19") "3.584"
Set-TextCell $ws.Range("E20") "-0.73%"
Set-TextCell $ws.Range("E21") "0.11%"
Set-TextCell $ws.Range("D22") "5.012"
Set-TextCell $ws.Range("E22") "2.37%"
Set-TextCell $ws.Range("D24") "0.04316"
Set-TextCell $ws.Range("E24") "-0.06%"
Set-TextCell $ws.Range("E25") "-0.37%"
Set-TextCell $ws.Range("D26") "0.004548"
Set-TextCell $ws.Range("E26") "-4.59%"
Set-TextCell $ws.Range("D27") "0.0001353"
Set-TextCell $ws.Range("E27") "4.19%"
Set-TextCell $ws.Range("D39") "0.02274"
Set-TextCell $ws.Range("E39") "2.39%"
Set-TextCell $ws.Range("E40") "-3.94%"
Set-TextCell $ws.Range("D41") "0.007447"
Set-TextCell $ws.Range("E41") "-0.93%"
Set-TextCell $ws.Range("D42") "0.009929"
Set-TextCell $ws.Range("E42") "-0.34%"
Set-TextCell $ws.Range("D43") "0.1354"
Set-TextCell $ws.Range("E43") "-1.82%"
Set-TextCell $ws.Range("D44") "0.002044"
Set-TextCell $ws.Range("E44") "-4.81%"
Set-TextCell $ws.Range("D45") "0.008769"
Set-TextCell $ws.Range("E45") "-10.97%"
Set-TextCell $ws.Range("D46") "0.00006762"
Set-TextCell $ws.Range("E46") "1.70%"
Set-TextCell $ws.Range("D48") "0.003003"
Set-TextCell $ws.Range("E48") "8.68%"
Set-TextCell $ws.Range("E49") "8.49%"
Set-TextCell $ws.Range("D50") "0.00002102"
Set-TextCell $ws.Range("D51") "0.0002002"
